$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9987067580223083
$ws.Range("B1").Value = 2.110029220581055
$ws.Range("C1").Value = 6.957849502563477
$ws.Range("D1").Value = 2.063240051269531
$ws.Range("E1").Value = 1.375236511230469
